$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (AD1:AF1), copying the existing header style from AC1
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins / Losses / Ties) for every player row
$wins = 83
$losses = 79
$ties = 0

for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
